$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'246.41"
$ws.Range("D3").Value = "'22.79"
$ws.Range("D4").Value = "'5.469"
$ws.Range("D5").Value = "'0.05744"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "'6.329"
$ws.Range("E7").Value = "6KuCoinTokenKCS"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.8124"
$ws.Range("E8").Value = "7MXTokenMX"
$ws.Range("B9").Value = "FTXToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D9").Value = "'0.8860"
$ws.Range("E9").Value = "8FTXTokenFTT"
$ws.Range("B10").Value = "WazirX"
$ws.Range("C10").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D10").Value = "'0.1441"
$ws.Range("E10").Value = "9WazirXWRX"
$ws.Range("B11").Value = "MandalaExchangeToken"
$ws.Range("C11").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D11").Value = "'0.07363"
$ws.Range("E11").Value = "10MandalaExchangeTokenMDX"
$ws.Range("D13").Value = "'0.03092"
$ws.Range("D15").Value = "'0.001590"
$ws.Range("D16").Value = "'0.04820"
$ws.Range("D17").Value = "'0.0005849"
$ws.Range("D19").Value = "'0.005117"
$ws.Range("D20").Value = "'0.0009964"
$ws.Range("D22").Value = "'3.751"
$ws.Range("B23").Value = "BTSEToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D23").Value = "'2.199"
$ws.Range("E23").Value = "22BTSETokenBTSE"
$ws.Range("B24").Value = "BitpandaEcosystemToken"
$ws.Range("C24").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D24").Value = "'0.3275"
$ws.Range("E24").Value = "23BitpandaEcosystemTokenBEST"
$ws.Range("B25").Value = "ProBitToken"
$ws.Range("C25").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D25").Value = "'0.1320"
$ws.Range("E25").Value = "24ProBitTokenPROB"
$ws.Range("D26").Value = "'4.179"
$ws.Range("D27").Value = "'0.0003158"
$ws.Range("D40").Value = "'0.03913"
$ws.Range("D41").Value = "'0.006756"
$ws.Range("E41").Value = "40KickTokenKICKBestin24h"
$ws.Range("D42").Value = "'0.1071"
$ws.Range("D43").Value = "'0.003199"
$ws.Range("E43").Value = "42CEJICEJI"
$ws.Range("D44").Value = "'0.007776"
$ws.Range("D45").Value = "'0.00005637"
$ws.Range("D47").Value = "'0.3799"
